$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "302.24" or
# thousand-dotted values like "42.852.36"). Force a text number format
# before assigning so Excel keeps the exact original digits/format instead
# of re-parsing the string into a floating point number, then restore the
# default "Normal" style so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.852.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.25%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.314.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.88%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.73%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.506"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.26%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -0.82%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.72%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.02"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0783"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.22%  "

$ws.Range("E13").Value = "  +0.11%  "

$ws.Range("E14").Value = "  -1.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.677.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.96%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.318.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.08%  "

$ws.Range("E17").Value = "  +2.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.788.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.03%  "

$ws.Range("E20").Value = "  +2.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0892"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.91%  "

$ws.Range("E23").Value = "  +6.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.37%  "

$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("E26").Value = "  +1.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.65%  "

$ws.Range("E28").Value = "  +14.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.51%  "

$ws.Range("E30").Value = "  +1.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.07%  "

$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("E33").Value = "  +0.99%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0699"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.19%  "

$ws.Range("E37").Value = "  -0.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.100"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("E40").Value = "  +1.42%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +14.45%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.927.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.32%  "

$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.48%  "

$ws.Range("E46").Value = "  -1.31%  "

$ws.Range("E47").Value = "  -0.40%  "

$ws.Range("E48").Value = "  -0.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.546.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.05%  "
